$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values to the worksheet (existing cells remain unchanged)
$ws.Range("H6").Value = "dfghgfd"
$ws.Range("H9").Value = "fghgf"
$ws.Range("F12").Value = "try"
$ws.Range("F13").Value = "xg"

# Update the selected cell to match the final state of the workbook
$ws.Range("F13").Select()
